$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 10 from 45185 (2023-09-16)
# to 45204 (2023-10-05), keeping the existing date formatting/style intact.
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45204
}
